# DLAD Part 25 edit: add the P25_802_71 bookmark around the "25.802-71"
# heading text (the runtime auto-renumbers every later bookmark id by +1,
# which reproduces the id cascade seen in the target revision).
#
# (The remaining hunks in the source revision are purely cosmetic
# byproducts of Word's background grammar/spell-check and repagination
# passes - proofErr wrappers and lastRenderedPageBreak markers around
# text that is byte-for-byte identical before/after - and carry no
# content change of their own.)

$d = $word.ActiveDocument

$target = "25.802-71 End use certificates."
$rng = $d.Content
$found = $rng.Find.Execute($target, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    Write-Output "ERROR: could not find heading text"
} else {
    $headingStart = $rng.Start
    Write-Output "Found '$($rng.Text)' at $headingStart"

    $bookmarkRange = $d.Range($headingStart, $headingStart + 9)
    Write-Output "Bookmark range text: '$($bookmarkRange.Text)'"

    $d.Bookmarks.Add("P25_802_71", $bookmarkRange)
    Write-Output "Bookmark P25_802_71 added; bookmark count now $($d.Bookmarks.Count)"
}
